# Add a new data row (row 5) to the Espinaca worksheet, mirroring the
# existing rows' layout. Row 3 is the most recent prior observation
# (2022-08-02 / 44775), this adds a newer weekly observation
# (2022-08-09 / 44782).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C5").Value = "Ñuble"
$ws.Range("D5").Value = 44782
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 100112012
$ws.Range("G5").Value = "Espinaca"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = 8500
$ws.Range("N5").Value = "$/cuna 10 kilos"
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 850
$ws.Range("Q5").Value = 10
$ws.Range("R5").Value = "Hortaliza"
